# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
# Re-orders/swaps match rows on the "Greece Super League 1" sheet.
# Column A (id) and row position stay fixed; columns B:AD (match data)
# are rearranged between rows as described below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B through AD (28 columns) hold the per-match data that moves.
$firstCol = 2   # B
$lastCol  = 30  # AD

function Get-RowData($row) {
    $data = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $data[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $data
}

function Set-RowData($row, $data) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $data[$c]
    }
}

# --- Swap row 124 <-> row 125 ---
$r124 = Get-RowData 124
$r125 = Get-RowData 125
Set-RowData 124 $r125
Set-RowData 125 $r124

# --- Swap row 170 <-> row 171 ---
$r170 = Get-RowData 170
$r171 = Get-RowData 171
Set-RowData 170 $r171
Set-RowData 171 $r170

# --- Cyclic rotation among rows 175..180 ---
# new175 = old176, new176 = old177, new177 = old178,
# new178 = old179, new179 = old180, new180 = old175
$r175 = Get-RowData 175
$r176 = Get-RowData 176
$r177 = Get-RowData 177
$r178 = Get-RowData 178
$r179 = Get-RowData 179
$r180 = Get-RowData 180

Set-RowData 175 $r176
Set-RowData 176 $r177
Set-RowData 177 $r178
Set-RowData 178 $r179
Set-RowData 179 $r180
Set-RowData 180 $r175

# --- Swap row 194 <-> row 195 ---
$r194 = Get-RowData 194
$r195 = Get-RowData 195
Set-RowData 194 $r195
Set-RowData 195 $r194
